# Update cryptos list (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to stay a text value (matches original inlineStr cells)
    # instead of being auto-coerced into a number by Excel's type inference,
    # while keeping the cell style back to the default (no explicit style).
    $range.NumberFormat = "@"
    $range.Value2 = $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "26.734.90"
Set-TextValue $ws.Range("E2") "  +1.18%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "1.724.95"
Set-TextValue $ws.Range("E3") "  +0.03%  "

# Row 4 - TetherUSD
Set-TextValue $ws.Range("D4") "0.9977"
Set-TextValue $ws.Range("E4") "  -0.20%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "240.69"
Set-TextValue $ws.Range("E5") "  -1.26%  "

# Row 6 - USDC
Set-TextValue $ws.Range("D6") "0.9980"
Set-TextValue $ws.Range("E6") "  -0.20%  "

# Row 7 - XRP
Set-TextValue $ws.Range("E7") "  -2.08%  "

# Row 8 - Cardano
Set-TextValue $ws.Range("D8") "0.2580"
Set-TextValue $ws.Range("E8") "  -1.44%  "

# Row 9 - Dogecoin
Set-TextValue $ws.Range("D9") "0.06194"
Set-TextValue $ws.Range("E9") "  -0.14%  "

# Row 10 - WrappedEther
Set-TextValue $ws.Range("D10") "1.726.23"
Set-TextValue $ws.Range("E10") "  +0.06%  "

# Row 11 - Solana
Set-TextValue $ws.Range("D11") "15.94"
Set-TextValue $ws.Range("E11") "  +3.12%  "

# Row 12 - TRON
Set-TextValue $ws.Range("D12") "0.06902"
Set-TextValue $ws.Range("E12") "  -1.33%  "

# Row 13 - Polygon
Set-TextValue $ws.Range("D13") "0.6066"
Set-TextValue $ws.Range("E13") "  +1.00%  "

# Row 14 - Polkadot
Set-TextValue $ws.Range("D14") "4.471"
Set-TextValue $ws.Range("E14") "  -1.70%  "

# Row 15 - Litecoin
Set-TextValue $ws.Range("D15") "76.86"
Set-TextValue $ws.Range("E15") "  -0.89%  "

# Row 16 - Dai
Set-TextValue $ws.Range("D16") "0.9983"
Set-TextValue $ws.Range("E16") "  -0.16%  "

# Row 17 - WrappedBTC
Set-TextValue $ws.Range("D17") "26.551.61"
Set-TextValue $ws.Range("E17") "  +0.45%  "

# Row 18 - BinanceUSD
Set-TextValue $ws.Range("D18") "0.9973"
Set-TextValue $ws.Range("E18") "  -0.25%  "

# Row 19 - ShibaInu
Set-TextValue $ws.Range("D19") "0.000007144"
Set-TextValue $ws.Range("E19") "  -0.78%  "

# Row 20 - Avalanche
Set-TextValue $ws.Range("E20") "  +0.60%  "

# Row 21 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D21") "1.949.10"
Set-TextValue $ws.Range("E21") "  +0.03%  "

# Row 22 - Uniswap
Set-TextValue $ws.Range("D22") "4.425"
Set-TextValue $ws.Range("E22") "  -1.07%  "

# Row 23 - Cosmos
Set-TextValue $ws.Range("D23") "8.566"
Set-TextValue $ws.Range("E23") "  -0.14%  "

# Row 24 - Chainlink
Set-TextValue $ws.Range("D24") "5.061"
Set-TextValue $ws.Range("E24") "  -1.99%  "

# Row 25 - Monero
Set-TextValue $ws.Range("D25") "137.13"
Set-TextValue $ws.Range("E25") "  -0.46%  "

# Row 26 - EthereumClassic
Set-TextValue $ws.Range("E26") "  -0.28%  "

# Row 27 - LidoDAOToken
Set-TextValue $ws.Range("D27") "1.768"
Set-TextValue $ws.Range("E27") "  +2.72%  "

# Row 28/29 - Toncoin and BitcoinCash swap places
Set-TextValue $ws.Range("B28") "Toncoin"
Set-TextValue $ws.Range("C28") "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D28") "1.380"
Set-TextValue $ws.Range("E28") "  -1.19%  "

Set-TextValue $ws.Range("B29") "BitcoinCash"
Set-TextValue $ws.Range("C29") "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws.Range("D29") "106.05"
Set-TextValue $ws.Range("E29") "  -0.98%  "

# Row 30 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("E30") "  -0.32%  "

# Row 31 - Stellar
Set-TextValue $ws.Range("D31") "0.07944"
Set-TextValue $ws.Range("E31") "  -0.86%  "

# Row 32 - Filecoin
Set-TextValue $ws.Range("D32") "3.690"
Set-TextValue $ws.Range("E32") "  +0.28%  "

# Row 33 - Hedera
Set-TextValue $ws.Range("D33") "0.04475"
Set-TextValue $ws.Range("E33") "  -1.07%  "

# Row 34 - Frax
Set-TextValue $ws.Range("D34") "0.9971"
Set-TextValue $ws.Range("E34") "  -0.20%  "

# Row 35 - HuobiToken
Set-TextValue $ws.Range("D35") "2.594"
Set-TextValue $ws.Range("E35") "  -0.29%  "

# Row 36 - ARBITRUM
Set-TextValue $ws.Range("D36") "1.005"
Set-TextValue $ws.Range("E36") "  +0.60%  "

# Row 37 - ImmutableX
Set-TextValue $ws.Range("D37") "0.6203"
Set-TextValue $ws.Range("E37") "  -1.08%  "

# Row 38 - TrustWalletToken
Set-TextValue $ws.Range("D38") "0.9279"
Set-TextValue $ws.Range("E38") "  -2.80%  "

# Row 39 - RenderToken
Set-TextValue $ws.Range("D39") "2.027"
Set-TextValue $ws.Range("E39") "  +4.07%  "

# Row 40 - MXToken
Set-TextValue $ws.Range("D40") "2.444"
Set-TextValue $ws.Range("E40") "  +2.20%  "

# Row 41 - PaxDollar
Set-TextValue $ws.Range("D41") "0.9972"
Set-TextValue $ws.Range("E41") "  -0.27%  "

# Row 42 - VeChain
Set-TextValue $ws.Range("E42") "  +0.46%  "

# Row 43 - FraxShare
Set-TextValue $ws.Range("D43") "5.649"
Set-TextValue $ws.Range("E43") "  +5.86%  "

# Row 44 - Quant
Set-TextValue $ws.Range("D44") "99.65"
Set-TextValue $ws.Range("E44") "  +0.13%  "

# Row 45 - TheSandbox
Set-TextValue $ws.Range("D45") "0.3830"
Set-TextValue $ws.Range("E45") "  -0.74%  "

# Row 46 - Aptos
Set-TextValue $ws.Range("D46") "6.840"
Set-TextValue $ws.Range("E46") "  +0.30%  "

# Row 47 - Algorand
Set-TextValue $ws.Range("D47") "0.1155"
Set-TextValue $ws.Range("E47") "  -1.20%  "

# Row 48 - Cronos
Set-TextValue $ws.Range("D48") "0.05390"

# Row 49 - EnergySwap
Set-TextValue $ws.Range("D49") "7.894"
Set-TextValue $ws.Range("E49") "  +1.87%  "

# Row 50 - Elrond
Set-TextValue $ws.Range("D50") "30.07"
Set-TextValue $ws.Range("E50") "  -0.58%  "

# Row 51 - Aave
Set-TextValue $ws.Range("D51") "51.45"
Set-TextValue $ws.Range("E51") "  +0.87%  "
